# Update the "Ativacao" date, and several syllabus / method description
# cells for LOB1257, per the 2022-01-09 site rebuild.
#
# Column B holds the "current" text, column C holds the "modified" (shown
# in red) text. In this workbook both columns already contain identical
# text for every changed row, so both get the same new value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Ativacao: 01/01/2020 -> 01/01/2022 --------------------------------
# Plain assignment would be auto-recognised as a date and stored as a
# serial number; round-trip it through a text-producing formula + paste
# values so it lands back as literal text (matches the original shared
# string type) instead of a date serial.
$dateCell = $ws.Range("Z1")
$dateCell.Formula = '="01/01/2022"'
$dateCell.Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$dateCell.ClearContents()
$excel.CutCopyMode = $false

# --- Programa resumido: / Short syllabus: -------------------------------
$ws.Range("B14").Value = 'Sistemas de Abastecimento de Água; Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água.'
$ws.Range("C14").Value = 'Sistemas de Abastecimento de Água; Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água.'

$ws.Range("B15").Value = 'Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.'
$ws.Range("C15").Value = 'Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations.'

# --- Programa: / Syllabus: ----------------------------------------------
$ws.Range("B16").Value = '- Tipos de dimensionamento de redes de distribuição de água;- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; - Padrão de Potabilidade; - Tecnologias de Tratamento de Água;- Unidades Constituintes de um Sistema de Abastecimento de Água;- Captação de Águas Subterrâneas e Captação de Águas Superficiais - Gradeamento, remoção de areia, casa de bombas; - Reservação; - Redes de distribuição: tipos de rede, peças e órgãos acessórios; - Sistema de Tratamento de Água de Ciclo Completo; - Coagulação-floculação e Mistura Rápida; - Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, - Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, - Tratamento de resíduos gerados nas ETAs e reuso de água recuperada'
$ws.Range("C16").Value = '- Tipos de dimensionamento de redes de distribuição de água;- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; - Padrão de Potabilidade; - Tecnologias de Tratamento de Água;- Unidades Constituintes de um Sistema de Abastecimento de Água;- Captação de Águas Subterrâneas e Captação de Águas Superficiais - Gradeamento, remoção de areia, casa de bombas; - Reservação; - Redes de distribuição: tipos de rede, peças e órgãos acessórios; - Sistema de Tratamento de Água de Ciclo Completo; - Coagulação-floculação e Mistura Rápida; - Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, - Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, - Tratamento de resíduos gerados nas ETAs e reuso de água recuperada'

$ws.Range("B17").Value = '- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water'
$ws.Range("C17").Value = '- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water'

# --- Metodo: ---------------------------------------------------------------
$ws.Range("B19").Value = 'Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios.'
$ws.Range("C19").Value = 'Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios.'
